# Updates the cryptos price/volume table (Sheet1) to the latest scraped
# snapshot. Only cells B2:E51 that actually changed values are touched;
# everything else (headers, row-index column A, styles) is left as-is.
#
# Price values in column D that look like plain decimal numbers
# (e.g. "6.69") are written with a leading apostrophe so Excel keeps
# them as text (matching the source data, which mixes thousands-dot
# formatted numbers like "60.899.61" with plain ones like "6.69" -- all
# stored as text). The apostrophe prefix nudges the cell into Excel's
# "quote prefix" text entry, so the style is then reset back to Normal
# to avoid leaving a stray Text number-format behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.899.61"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").Value = "2.599.19"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'522.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.38%  "
$ws.Range("D6").Value = "'154.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.93%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +2.23%  "
$ws.Range("D9").Value = "'6.69"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.12%  "
$ws.Range("E10").Value = "  +2.51%  "
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("E12").Value = "  +1.57%  "
$ws.Range("D13").Value = "3.055.32"
$ws.Range("E13").Value = "  +0.70%  "
$ws.Range("D14").Value = "60.934.43"
$ws.Range("E14").Value = "  +1.28%  "
$ws.Range("D15").Value = "'21.71"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.26%  "
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("D17").Value = "2.608.07"
$ws.Range("E17").Value = "  +0.83%  "
$ws.Range("E18").Value = "  -1.31%  "
$ws.Range("D19").Value = "'352.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.21%  "
$ws.Range("E20").Value = "  +1.93%  "
$ws.Range("D21").Value = "'6.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").Value = "'61.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.73%  "
$ws.Range("D24").Value = "'0.427"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.64%  "
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("D26").Value = "2.716.95"
$ws.Range("E26").Value = "  +0.64%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("D28").Value = "0.0₃0848"
$ws.Range("E28").Value = "  +0.82%  "
$ws.Range("E29").Value = "  +0.84%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").Value = "'6.34"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.57%  "
$ws.Range("D32").Value = "'19.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.58%  "
$ws.Range("D33").Value = "'1.60"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.29%  "
$ws.Range("D34").Value = "'149.36"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.67%  "
$ws.Range("E35").Value = "  +6.56%  "
$ws.Range("D36").Value = "'0.948"
$ws.Range("D36").Style = "Normal"
$ws.Range("E37").Value = "  +1.62%  "
$ws.Range("E38").Value = "  +2.39%  "
$ws.Range("D39").Value = "'0.850"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.40%  "
$ws.Range("D40").Value = "'3.79"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.01%  "
$ws.Range("E41").Value = "  +1.70%  "
$ws.Range("D42").Value = "'287.43"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.25%  "
$ws.Range("E43").Value = "  +1.55%  "
$ws.Range("E44").Value = "  +2.04%  "
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("D46").Value = "'0.998"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'19.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.52%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'4.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("D49").Value = "'0.0237"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.93%  "
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("D51").Value = "'19.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.39%  "
